$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.71%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.800"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.27%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.46%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.957"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.58%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.275"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-11.69%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8765"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.02%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1540"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.03%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05056"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.65%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07517"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.57%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03018"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.76%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09038"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.21%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.46%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006402"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.19%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005861"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.52%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.453"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.99%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.299"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.61%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.30%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.947"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.81%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.19%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001171"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.60%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003864"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.97%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001935"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "19.65%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04155"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.04%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006849"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.87%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.42%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-14.56%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01114"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.92%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005168"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.93%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02298"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "8.28%"
